$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.040175333333333
$ws.Range("H2").Value = 15.120526
$ws.Range("I2").Value = 0.1638830352839606
$ws.Range("J2").Value = 0.1638830352839606
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.180428
$ws.Range("N2").Value = 0.541284
$ws.Range("O2").Value = 0.6724860231084607
$ws.Range("P2").Value = 0.6724860231084607
$ws.Range("Q2").Value = 0.9093887550426666
$ws.Range("R2").Value = 8.184498795384
$ws.Range("S2").Value = 0.1102090506530542
$ws.Range("T2").Value = 0.1102090506530542

$ws.Range("G3").Value = 5.040175333333333
$ws.Range("H3").Value = 15.120526
$ws.Range("I3").Value = 0.1638830352839606
$ws.Range("J3").Value = 0.1638830352839606
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01727566666666666
$ws.Range("N3").Value = 0.051827
$ws.Range("O3").Value = 0.06438936513852653
$ws.Range("P3").Value = 0.06438936513852653
$ws.Range("Q3").Value = 0.08707238900022221
$ws.Range("R3").Value = 0.783651501002
$ws.Range("S3").Value = 0.01055232459890897
$ws.Range("T3").Value = 0.01055232459890897

$ws.Range("G4").Value = 5.040175333333333
$ws.Range("H4").Value = 15.120526
$ws.Range("I4").Value = 0.1638830352839606
$ws.Range("J4").Value = 0.1638830352839606
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07059633333333333
$ws.Range("N4").Value = 0.211789
$ws.Range("O4").Value = 0.2631246117530128
$ws.Range("P4").Value = 0.2631246117530128
$ws.Range("Q4").Value = 0.3558178978904444
$ws.Range("R4").Value = 3.202361081014
$ws.Range("S4").Value = 0.04312166003199744
$ws.Range("T4").Value = 0.04312166003199745

$ws.Range("G5").Value = 6.048190666666667
$ws.Range("H5").Value = 18.144572
$ws.Range("I5").Value = 0.1966590007046292
$ws.Range("J5").Value = 0.1966590007046292
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.180428
$ws.Range("N5").Value = 0.541284
$ws.Range("O5").Value = 0.6724860231084607
$ws.Range("P5").Value = 0.6724860231084607
$ws.Range("Q5").Value = 1.091262945605333
$ws.Range("R5").Value = 9.821366510448
$ws.Range("S5").Value = 0.1322504292923401
$ws.Range("T5").Value = 0.1322504292923401

$ws.Range("G6").Value = 6.048190666666667
$ws.Range("H6").Value = 18.144572
$ws.Range("I6").Value = 0.1966590007046292
$ws.Range("J6").Value = 0.1966590007046292
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01727566666666666
$ws.Range("N6").Value = 0.051827
$ws.Range("O6").Value = 0.06438936513852653
$ws.Range("P6").Value = 0.06438936513852653
$ws.Range("Q6").Value = 0.1044865258937778
$ws.Range("R6").Value = 0.940378733044
$ws.Range("S6").Value = 0.01266274820414811
$ws.Range("T6").Value = 0.01266274820414812

$ws.Range("G7").Value = 6.048190666666667
$ws.Range("H7").Value = 18.144572
$ws.Range("I7").Value = 0.1966590007046292
$ws.Range("J7").Value = 0.1966590007046292
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.07059633333333333
$ws.Range("N7").Value = 0.211789
$ws.Range("O7").Value = 0.2631246117530128
$ws.Range("P7").Value = 0.2631246117530128
$ws.Range("Q7").Value = 0.4269800843675555
$ws.Range("R7").Value = 3.842820759308
$ws.Range("S7").Value = 0.05174582320814103
$ws.Range("T7").Value = 0.05174582320814104

$ws.Range("G8").Value = 4.433369666666667
$ws.Range("H8").Value = 13.300109
$ws.Range("I8").Value = 0.1441525402309101
$ws.Range("J8").Value = 0.1441525402309101
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.180428
$ws.Range("N8").Value = 0.541284
$ws.Range("O8").Value = 0.6724860231084607
$ws.Range("P8").Value = 0.6724860231084607
$ws.Range("Q8").Value = 0.7999040222173334
$ws.Range("R8").Value = 7.199136199956
$ws.Range("S8").Value = 0.09694056850086714
$ws.Range("T8").Value = 0.09694056850086714

$ws.Range("G9").Value = 4.433369666666667
$ws.Range("H9").Value = 13.300109
$ws.Range("I9").Value = 0.1441525402309101
$ws.Range("J9").Value = 0.1441525402309101
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01727566666666666
$ws.Range("N9").Value = 0.051827
$ws.Range("O9").Value = 0.06438936513852653
$ws.Range("P9").Value = 0.06438936513852653
$ws.Range("Q9").Value = 0.07658941657144444
$ws.Range("R9").Value = 0.6893047491429999
$ws.Range("S9").Value = 0.009281890548574206
$ws.Range("T9").Value = 0.009281890548574206

$ws.Range("G10").Value = 4.433369666666667
$ws.Range("H10").Value = 13.300109
$ws.Range("I10").Value = 0.1441525402309101
$ws.Range("J10").Value = 0.1441525402309101
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.07059633333333333
$ws.Range("N10").Value = 0.211789
$ws.Range("O10").Value = 0.2631246117530128
$ws.Range("P10").Value = 0.2631246117530128
$ws.Range("Q10").Value = 0.3129796427778889
$ws.Range("R10").Value = 2.816816785001
$ws.Range("S10").Value = 0.03793008118146879
$ws.Range("T10").Value = 0.03793008118146879

$ws.Range("G11").Value = 3.607224333333333
$ws.Range("H11").Value = 10.821673
$ws.Range("I11").Value = 0.1172901404415748
$ws.Range("J11").Value = 0.1172901404415748
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.180428
$ws.Range("N11").Value = 0.541284
$ws.Range("O11").Value = 0.6724860231084607
$ws.Range("P11").Value = 0.6724860231084607
$ws.Range("Q11").Value = 0.6508442720146667
$ws.Range("R11").Value = 5.857598448131999
$ws.Range("S11").Value = 0.0788759800953875
$ws.Range("T11").Value = 0.0788759800953875

$ws.Range("G12").Value = 3.607224333333333
$ws.Range("H12").Value = 10.821673
$ws.Range("I12").Value = 0.1172901404415748
$ws.Range("J12").Value = 0.1172901404415748
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.01727566666666666
$ws.Range("N12").Value = 0.051827
$ws.Range("O12").Value = 0.06438936513852653
$ws.Range("P12").Value = 0.06438936513852653
$ws.Range("Q12").Value = 0.06231720517455554
$ws.Range("R12").Value = 0.5608548465709999
$ws.Range("S12").Value = 0.00755223768004162
$ws.Range("T12").Value = 0.00755223768004162

$ws.Range("G13").Value = 3.607224333333333
$ws.Range("H13").Value = 10.821673
$ws.Range("I13").Value = 0.1172901404415748
$ws.Range("J13").Value = 0.1172901404415748
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.07059633333333333
$ws.Range("N13").Value = 0.211789
$ws.Range("O13").Value = 0.2631246117530128
$ws.Range("P13").Value = 0.2631246117530128
$ws.Range("Q13").Value = 0.2546568114441111
$ws.Range("R13").Value = 2.291911302997
$ws.Range("S13").Value = 0.03086192266614573
$ws.Range("T13").Value = 0.03086192266614573

$ws.Range("G14").Value = 6.973136333333334
$ws.Range("H14").Value = 20.919409
$ws.Range("I14").Value = 0.2267339273294199
$ws.Range("J14").Value = 0.2267339273294199
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.180428
$ws.Range("N14").Value = 0.541284
$ws.Range("O14").Value = 0.6724860231084607
$ws.Range("P14").Value = 0.6724860231084607
$ws.Range("Q14").Value = 1.258149042350667
$ws.Range("R14").Value = 11.323341381156
$ws.Range("S14").Value = 0.1524753970935243
$ws.Range("T14").Value = 0.1524753970935243

$ws.Range("G15").Value = 6.973136333333334
$ws.Range("H15").Value = 20.919409
$ws.Range("I15").Value = 0.2267339273294199
$ws.Range("J15").Value = 0.2267339273294199
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.3333333333333333
$ws.Range("M15").Value = 0.01727566666666666
$ws.Range("N15").Value = 0.051827
$ws.Range("O15").Value = 0.06438936513852653
$ws.Range("P15").Value = 0.06438936513852653
$ws.Range("Q15").Value = 0.1204655789158889
$ws.Range("R15").Value = 1.084190210243
$ws.Range("S15").Value = 0.01459925363610616
$ws.Range("T15").Value = 0.01459925363610616

$ws.Range("G16").Value = 6.973136333333334
$ws.Range("H16").Value = 20.919409
$ws.Range("I16").Value = 0.2267339273294199
$ws.Range("J16").Value = 0.2267339273294199
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.07059633333333333
$ws.Range("N16").Value = 0.211789
$ws.Range("O16").Value = 0.2631246117530128
$ws.Range("P16").Value = 0.2631246117530128
$ws.Range("Q16").Value = 0.4922778569667778
$ws.Range("R16").Value = 4.430500712701001
$ws.Range("S16").Value = 0.05965927659978943
$ws.Range("T16").Value = 0.05965927659978943

$ws.Range("G17").Value = 4.652614333333333
$ws.Range("H17").Value = 13.957843
$ws.Range("I17").Value = 0.1512813560095054
$ws.Range("J17").Value = 0.1512813560095054
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.180428
$ws.Range("N17").Value = 0.541284
$ws.Range("O17").Value = 0.6724860231084607
$ws.Range("P17").Value = 0.6724860231084607
$ws.Range("Q17").Value = 0.8394618989346667
$ws.Range("R17").Value = 7.555157090412
$ws.Range("S17").Value = 0.1017345974732875
$ws.Range("T17").Value = 0.1017345974732875

$ws.Range("G18").Value = 4.652614333333333
$ws.Range("H18").Value = 13.957843
$ws.Range("I18").Value = 0.1512813560095054
$ws.Range("J18").Value = 0.1512813560095054
$ws.Range("K18").Value = 1
$ws.Range("L18").Value = 0.3333333333333333
$ws.Range("M18").Value = 0.01727566666666666
$ws.Range("N18").Value = 0.051827
$ws.Range("O18").Value = 0.06438936513852653
$ws.Range("P18").Value = 0.06438936513852653
$ws.Range("Q18").Value = 0.08037701435122221
$ws.Range("R18").Value = 0.723393129161
$ws.Range("S18").Value = 0.009740910470747467
$ws.Range("T18").Value = 0.009740910470747469

$ws.Range("G19").Value = 4.652614333333333
$ws.Range("H19").Value = 13.957843
$ws.Range("I19").Value = 0.1512813560095054
$ws.Range("J19").Value = 0.1512813560095054
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.07059633333333333
$ws.Range("N19").Value = 0.211789
$ws.Range("O19").Value = 0.2631246117530128
$ws.Range("P19").Value = 0.2631246117530128
$ws.Range("Q19").Value = 0.3284575123474444
$ws.Range("R19").Value = 2.956117611127
$ws.Range("S19").Value = 0.03980584806547043
$ws.Range("T19").Value = 0.03980584806547043
